# Weekly price-list update: a new daily record is inserted at the top of
# the data (row 23, just after the header block that stays put through
# row 22), pushing every existing record down by one row. The new record
# carries the latest date (2022-01-04 -> serial 44565) and its own
# min/max/avg/ $-per-kg figures; all other fields mirror the record that
# used to occupy row 23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 23:96 down to 24:97, leaving a blank row 23 to fill in.
$ws.Rows.Item(23).Insert()

$ws.Range("A23").Value = 11
$ws.Range("B23").Value = "Vega Monumental Concepción"
$ws.Range("C23").Value = "Bíobío"
$ws.Range("D23").Value = 44565
$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 100112043
$ws.Range("G23").Value = "Pepino ensalada"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 7000
$ws.Range("L23").Value = 8000
$ws.Range("M23").Value = 7500
$ws.Range("N23").Value = "`$/caja 60 unidades"
$ws.Range("O23").Value = "Región de Arica y Parinacota"
$ws.Range("P23").Value = 125
$ws.Range("Q23").Value = 60
$ws.Range("R23").Value = "Hortaliza"
